$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (rows 2-7) with 1, matching the new "E" values added next
# to the existing cross/mean/SE/sex/environ/sire/dam table.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 5).Value = 1
}

# Update the active selection to E8, as recorded in the saved sheet view.
$ws.Range("E8").Select()
